$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C22").Value = "Öğrenci Tablosunda non clustered not unique index, Fakülte tablosunda nonclustered unique index kullanımı vardır. Ogretmen tablosunda da vardir."

$ws.Rows.Item(22).RowHeight = 30

$ws.Range("C22").Select()
